# GH-305: Extract class, add service information sheet
#
# Adds a new worksheet "Служебная информация" (Service information) at the
# end of the workbook, containing three label/value rows describing the
# R7.University export (version, export datetime, exporting user), and
# makes that new sheet the active one (mirrors the author's commit).

$wb = $excel.ActiveWorkbook

# --- Create the new sheet after the last existing one -----------------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$infoSheet = $wb.Worksheets.Add($null, $lastSheet)
$infoSheet.Name = "Служебная информация"

# --- Column widths (characters) matching the template ------------------
$infoSheet.Columns.Item(1).ColumnWidth = 36.95
$infoSheet.Columns.Item(2).ColumnWidth = 26.13

# --- Content: label column (A) / placeholder value column (B) ----------
$infoSheet.Range("A1").Value = "Версия R7.University:"
$infoSheet.Range("B1").Value = "{{UniversityVersion}}"

$infoSheet.Range("A2").Value = "Данные выгружены:"
$infoSheet.Range("B2").Value = "{{DataExportedAtDateTime}}"

$infoSheet.Range("A3").Value = "Пользователь, выгрузивший данные:"
$infoSheet.Range("B3").Value = "{{DataExportedByUserName}}"

# Labels in column A are bold, like the rest of the workbook's templates.
$infoSheet.Range("A1:A3").Font.Bold = $true

# --- Page setup, matching the other sheets in this template -------------
$ps = $infoSheet.PageSetup
$ps.LeftMargin = 56.7
$ps.RightMargin = 56.7
$ps.TopMargin = 75.8
$ps.BottomMargin = 75.8
$ps.HeaderMargin = 56.7
$ps.FooterMargin = 56.7
$ps.CenterHorizontally = $false
$ps.CenterVertically = $false
$ps.PrintGridlines = $false
$ps.PrintHeadings = $false
$ps.Orientation = 1
$ps.FitToPagesWide = 1
$ps.FitToPagesTall = 1
$ps.PaperSize = 9
$ps.CenterHeader = "&""Times New Roman,Обычный""&12&A"
$ps.CenterFooter = "&""Times New Roman,Обычный""&12Страница &P"

# --- Selection on the new sheet (matches the authored selection) -------
[void]$infoSheet.Range("B4").Select()

# The new sheet becomes the active tab, like in the commit.
$infoSheet.Activate()
